$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Range("C2").Value = "siva's id"
$ws.Range("A3").Value = "DCATEST4"
$ws.Range("B3").Value = "Password#1"
$ws.Range("C3").Value = "ramya's id"

$ws.Columns.Item(3).ColumnWidth = 14.28515625

$ws.Range("B6").Select()
